$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; everything currently at/after row 15
# (rows 15-26) shifts down to rows 16-27.
$ws.Rows(15).Insert()

# Populate the newly inserted row 15 with the new weekly price record.
$ws.Range("A15").Value = 7
$ws.Range("B15").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C15").Value = "Ñuble"
$ws.Range("D15").Value = 45280
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100102
$ws.Range("H15").Value = "Cítricos"
$ws.Range("I15").Value = 100102006
$ws.Range("J15").Value = "Pomelo"
$ws.Range("K15").Value = "Start Ruby"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("Q15").Value = "`$/caja 14 kilos empedrada"
$ws.Range("R15").Value = "Región de O'Higgins"
$ws.Range("S15").Value = 1071
$ws.Range("T15").Value = 14
